# Add a new sheet "Ambika Jewellers 004" (an invoice) at the end of the
# workbook, mirroring the layout/style of the existing "Shezwan House 003"
# sheet, and make it the active sheet/tab.

$wb  = $excel.ActiveWorkbook
$ws3 = $wb.Worksheets.Item(3)

# --- 1. Move the stored selection on "Shezwan House 003" before we change
#        the active sheet, so it ends up as a plain (non-active) sheetView
#        with its selection parked at B8 -- matching the target diff.
$ws3.Range("B8").Select() | Out-Null

# --- 2. Create the new worksheet as the last tab.
$count = $wb.Worksheets.Count
$lastSheet = $wb.Worksheets.Item($count)
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws.Name = "Ambika Jewellers 004"

# --- 3. Populate values/formulas FIRST (before copying any formatting in),
#        so the formula dependency graph is built cleanly.
$ws.Range("A1").Value = "SR NO"
$ws.Range("B1").Value = "ITEM DESCRIPTION"
$ws.Range("C1").Value = "QTY"
$ws.Range("D1").Value = "PRICE"
$ws.Range("E1").Value = "AMOUNT"

$ws.Range("A2").Value = 1
$ws.Range("B2").Value = "CCTV Cable 3+1"
$ws.Range("C2").Value = 10
$ws.Range("D2").Value = 90

$ws.Range("A3").Value = 2
$ws.Range("B3").Value = "Service Calls Charges"
$ws.Range("C3").Value = 1
$ws.Range("D3").Value = 900

$ws.Range("E2:E3").Formula = "=C2*D2"

$ws.Range("A4").Value = "TOTAL"
$ws.Range("E4").Formula = "=SUM(E2:E3)"

# --- 4. Copy formatting (styles only) from the equivalent rows of
#        "Shezwan House 003" onto the new sheet.
$ws3.Range("A1:E1").Copy() | Out-Null
$ws.Range("A1:E1").PasteSpecial(-4122) | Out-Null

$ws3.Range("A2:E3").Copy() | Out-Null
$ws.Range("A2:E3").PasteSpecial(-4122) | Out-Null

$ws3.Range("A10:E10").Copy() | Out-Null
$ws.Range("A4:E4").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = $false

# --- 5. Merge the TOTAL label row, like the source sheet.
$ws.Range("A4:D4").Merge() | Out-Null

# --- 6. Row height / column width cosmetics.
$ws.Range("A1").RowHeight = 25.8
$ws.Columns.Item(2).ColumnWidth = 21.7

# --- 7. Leave the new sheet selected/active with the whole table selected,
#        matching the target's <selection sqref="A1:E4"/>.
$ws.Range("A1:E4").Select() | Out-Null

Write-Output "Added sheet 'Ambika Jewellers 004'"
